$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "list"

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "firstName"
$ws.Range("B1").Value = "middleName"
$ws.Range("C1").Value = "lastName"
$ws.Range("D1").Value = "userName"
$ws.Range("E1").Value = "password"
$ws.Range("F1").ClearContents()

# ---- Row 2 ----
$ws.Range("A2").Value = "babar14"
$ws.Range("B2").Value = "babar26"
$ws.Range("C2").Value = "babar38"
$ws.Range("D2").Value = "abcxyz83"
$ws.Range("E2").Value = "aaAA1234eo05"
$ws.Range("F2").ClearContents()

# ---- Row 3 ----
$ws.Range("A3").Value = "pakin"
$ws.Range("B3").Value = "mser"
$ws.Range("C3").Value = "lsflsf"
$ws.Range("D3").Value = "sk4i"
$ws.Range("E3").Value = "lf34d`$4#`$"
$ws.Range("F3").ClearContents()

# ---- Row 4 ----
$ws.Range("A4").Value = "injg4"
$ws.Range("B4").Value = "sdf4rg"
$ws.Range("C4").Value = "wf344"
$ws.Range("D4").Value = "ret334"
$ws.Range("E4").Value = "5t4f4r5Frfg"
$ws.Range("F4").ClearContents()

# Apply the new font colour + vertical-center alignment to A3:F4.
# Build the combined format on a single cell first so only one new
# style record is produced, then fan it out with a format-only paste.
$seed = $ws.Range("A3")
$seed.Font.Color = 0
$seed.VerticalAlignment = -4108
$seed.Copy()
$ws.Range("A3:F4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths for E and F
$ws.Columns.Item(5).ColumnWidth = 12.3
$ws.Columns.Item(6).ColumnWidth = 16

# Selection
$ws.Range("E4").Select()
